# "Generate Report for Handoff"
# The localization status report moves from "In Translation" to
# "Ready for handoff" and the handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status columns: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest Handoff Datetime refreshed for both locales
$wsZhCn.Range("H2").Value = "2016-08-22 15:08:16"
$wsDeDe.Range("H2").Value = "2016-08-22 15:08:23"

# Overview's "Latest HO Xliff Generate Date" mirrors the de-de handoff datetime
$wsOverview.Range("G2").Value = "2016-08-22 15:08:23"

# Columns widened to fit the longer "Ready for handoff" status text
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
